# Scheduled market-data refresh: update currentAveragePrice / LevePrice /
# LeveProfit columns (H:N) on various Leve rows across the ALC, ARM, BSM,
# CRP, CUL, GSM, LTW and WVR sheets to reflect newly pulled prices.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 98 on ALC
$ws.Range("H98").Value = 3392.6086
$ws.Range("I98").Value = 3918.6316
$ws.Range("J98").Value = 894
$ws.Range("K98").Value = 3918.6316
$ws.Range("L98").Value = 894
$ws.Range("M98").Value = -2420.6316
$ws.Range("N98").Value = -3890

# Row 116 on ALC
$ws.Range("H116").Value = 2484
$ws.Range("I116").Value = 1675.8462
$ws.Range("J116").Value = 3984.8572
$ws.Range("K116").Value = 1675.8462
$ws.Range("L116").Value = 3984.8572
$ws.Range("M116").Value = 1766.1538
$ws.Range("N116").Value = -10868.8572

# Row 122 on ALC
$ws.Range("H122").Value = 3392.6086
$ws.Range("I122").Value = 3918.6316
$ws.Range("J122").Value = 894
$ws.Range("K122").Value = 11755.8948
$ws.Range("L122").Value = 2682
$ws.Range("M122").Value = -9305.8948
$ws.Range("N122").Value = -7582

# Row 132 on ALC
$ws.Range("H132").Value = 9016147
$ws.Range("I132").Value = 12351989
$ws.Range("J132").Value = 9372.4
$ws.Range("K132").Value = 37055967
$ws.Range("L132").Value = 28117.2
$ws.Range("M132").Value = -37053437
$ws.Range("N132").Value = -33177.2

# Row 138 on ALC
$ws.Range("H138").Value = 1721.9535
$ws.Range("I138").Value = 1974
$ws.Range("K138").Value = 5922
$ws.Range("M138").Value = -782

# Row 141 on ALC
$ws.Range("H141").Value = 4148
$ws.Range("I141").Value = 4506
$ws.Range("K141").Value = 13518
$ws.Range("M141").Value = -8338

$ws = $wb.Worksheets.Item("ARM")
# Row 32 on ARM
$ws.Range("H32").Value = 1544.04
$ws.Range("I32").Value = 1415.3118
$ws.Range("J32").Value = 3254.2856
$ws.Range("K32").Value = 1415.3118
$ws.Range("L32").Value = 3254.2856
$ws.Range("M32").Value = -1128.3118
$ws.Range("N32").Value = -3828.2856

# Row 63 on ARM
$ws.Range("H63").Value = 16950916
$ws.Range("I63").Value = 1604.1915
$ws.Range("J63").Value = 83335720
$ws.Range("K63").Value = 1604.1915
$ws.Range("L63").Value = 83335720
$ws.Range("M63").Value = -918.1914999999999
$ws.Range("N63").Value = -83337092

# Row 66 on ARM
$ws.Range("H66").Value = 16950916
$ws.Range("I66").Value = 1604.1915
$ws.Range("J66").Value = 83335720
$ws.Range("K66").Value = 8020.9575
$ws.Range("L66").Value = 416678600
$ws.Range("M66").Value = -4588.9575
$ws.Range("N66").Value = -416685464

# Row 122 on ARM
$ws.Range("H122").Value = 2179.2354
$ws.Range("I122").Value = 2070.4443
$ws.Range("J122").Value = 2301.625
$ws.Range("K122").Value = 6211.3329
$ws.Range("L122").Value = 6904.875
$ws.Range("M122").Value = -3761.3329
$ws.Range("N122").Value = -11804.875

# Row 132 on ARM
$ws.Range("H132").Value = 2054.739
$ws.Range("I132").Value = 1625.0646
$ws.Range("J132").Value = 2942.7334
$ws.Range("K132").Value = 4875.1938
$ws.Range("L132").Value = 8828.200199999999
$ws.Range("M132").Value = -2345.1938
$ws.Range("N132").Value = -13888.2002

# Row 133 on ARM
$ws.Range("H133").Value = 33420
$ws.Range("J133").Value = 33420
$ws.Range("L133").Value = 33420
$ws.Range("N133").Value = -38480

$ws = $wb.Worksheets.Item("BSM")
# Row 107 on BSM
$ws.Range("H107").Value = 1122.08
$ws.Range("I107").Value = 890.381
$ws.Range("J107").Value = 2338.5
$ws.Range("K107").Value = 890.381
$ws.Range("L107").Value = 2338.5
$ws.Range("M107").Value = 1029.619
$ws.Range("N107").Value = -6178.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31 on CRP
$ws.Range("H31").Value = 1462.6136
$ws.Range("J31").Value = 1921.3334
$ws.Range("L31").Value = 1921.3334
$ws.Range("N31").Value = -2511.3334

# Row 34 on CRP
$ws.Range("H34").Value = 1462.6136
$ws.Range("J34").Value = 1921.3334
$ws.Range("L34").Value = 1921.3334
$ws.Range("N34").Value = -2325.3334

$ws = $wb.Worksheets.Item("CUL")
# Row 113 on CUL
$ws.Range("H113").Value = 682.37933
$ws.Range("J113").Value = 726.15
$ws.Range("L113").Value = 2178.45
$ws.Range("N113").Value = -6518.45

# Row 131 on CUL
$ws.Range("H131").Value = 16950070
$ws.Range("J131").Value = 1076.8372
$ws.Range("L131").Value = 3230.5116
$ws.Range("N131").Value = -13310.5116

$ws = $wb.Worksheets.Item("GSM")
# Row 12 on GSM
$ws.Range("H12").Value = 6194647
$ws.Range("I12").Value = 6144312.5
$ws.Range("J12").Value = 7000000
$ws.Range("K12").Value = 6144312.5
$ws.Range("L12").Value = 7000000
$ws.Range("M12").Value = -6144172.5
$ws.Range("N12").Value = -7000280

# Row 19 on GSM
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

# Row 113 on GSM
$ws.Range("H113").Value = 1578.3334
$ws.Range("I113").Value = 1451.1111
$ws.Range("J113").Value = 1960
$ws.Range("K113").Value = 1451.1111
$ws.Range("L113").Value = 1960
$ws.Range("M113").Value = 718.8888999999999
$ws.Range("N113").Value = -6300

# Row 122 on GSM
$ws.Range("H122").Value = 3191.818
$ws.Range("I122").Value = 3557.077
$ws.Range("J122").Value = 2664.2222
$ws.Range("K122").Value = 10671.231
$ws.Range("L122").Value = 7992.6666
$ws.Range("M122").Value = -8221.231
$ws.Range("N122").Value = -12892.6666

$ws = $wb.Worksheets.Item("LTW")
# Row 7 on LTW
$ws.Range("H7").Value = 2068.7144
$ws.Range("J7").Value = 2397
$ws.Range("L7").Value = 2397
$ws.Range("N7").Value = -2621

# Row 22 on LTW
$ws.Range("H22").Value = 983.7273
$ws.Range("I22").Value = 855.25
$ws.Range("J22").Value = 1057.1428
$ws.Range("K22").Value = 855.25
$ws.Range("L22").Value = 1057.1428
$ws.Range("M22").Value = -560.25
$ws.Range("N22").Value = -1647.1428

# Row 27 on LTW
$ws.Range("H27").Value = 983.7273
$ws.Range("I27").Value = 855.25
$ws.Range("J27").Value = 1057.1428
$ws.Range("K27").Value = 855.25
$ws.Range("L27").Value = 1057.1428
$ws.Range("M27").Value = -748.25
$ws.Range("N27").Value = -1271.1428

# Row 93 on LTW
$ws.Range("H93").Value = 799
$ws.Range("I93").Value = 799
$ws.Range("K93").Value = 799
$ws.Range("M93").Value = 449

# Row 122 on LTW
$ws.Range("H122").Value = 83335336
$ws.Range("I122").Value = 83335336
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 250006008
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -250003558
$ws.Range("N122").ClearContents()

# Row 126 on LTW
$ws.Range("H126").Value = 2068.7144
$ws.Range("J126").Value = 2397
$ws.Range("L126").Value = 7191
$ws.Range("N126").Value = -12131

$ws = $wb.Worksheets.Item("WVR")
# Row 96 on WVR
$ws.Range("H96").Value = 1626.6
$ws.Range("I96").Value = 2022
$ws.Range("K96").Value = 2022
$ws.Range("M96").Value = -649

# Row 122 on WVR
$ws.Range("H122").Value = 19232602
$ws.Range("I122").Value = 22729084
$ws.Range("K122").Value = 68187252
$ws.Range("M122").Value = -68184802

# Row 126 on WVR
$ws.Range("H126").Value = 71433560
$ws.Range("I126").Value = 90912070
$ws.Range("J126").Value = 12333.333
$ws.Range("K126").Value = 272736210
$ws.Range("L126").Value = 36999.999
$ws.Range("M126").Value = -272733740
$ws.Range("N126").Value = -41939.999
